# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型"
# sheets to reflect the latest generated data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 313
$ws1.Range("F4").Value = 234
$ws1.Range("F5").Value = 2767
$ws1.Range("F6").Value = 1938
$ws1.Range("F8").Value = 123
$ws1.Range("F9").Value = 986
$ws1.Range("F10").Value = 185
$ws1.Range("F11").Value = 26

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 313
$ws4.Range("F4").Value = 234
$ws4.Range("F5").Value = 2767
$ws4.Range("F6").Value = 1938
$ws4.Range("F9").Value = 123
$ws4.Range("F10").Value = 986
$ws4.Range("F11").Value = 185
$ws4.Range("F12").Value = 26
